$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 19
$ws.Range("A19").Value = "Hazell et al 2020 WP (w/ Nakamura and Steinsson)"
$ws.Range("B19").Value = "Estimate slope of PC, show it's always been flat, even early 1980s. No missing disinflation or missing reinflation."
$ws.Range("C19").Value = "Results imply that drop in core inflation in 1980s due to shifting expectations about long-run mon pol as opposed to a steep Phillips curve, and greater stability of inflation since 1990s is mostly due to long-run infl expectations becoming more firmly anchored. "

# Row 20
$ws.Range("A20").Value = "Boutros et al 2020 WP"
$ws.Range("B20").Value = "Use forecasts of one-year S&P500 returns to track how beliefs of chief financial officers (CFOs) evolve. CFOs' beliefs are unbiased on average, but have too narrow Cis (they refer to this as miscalibration). When returns realized fall outside the Cis, the CIs widen."
$ws.Range("C20").Value = "Interpret this as consistent with Bayesian learning. But the magnitude of updating is too small. They say it is `"dampened by strong conviction in beliefs in the initial miscalibration`" and as a result, miscalibration persists. I.e. priors are really strong."

# Match style (wrap text) and row height of the other data rows
$ws.Range("A19:C20").WrapText = $true
$ws.Rows.Item(19).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 45

# Update the view to reflect the scrolled position / selection from the diff
# (new rows pushed the sheet so the window now shows row 13 at the top,
# with the next empty row, A21, selected)
$excel.Goto($ws.Range("A13"), $true)
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("A21").Select()
